$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Project Planner")

# Update cell values in the activity table
$ws.Range("F8").Value = 2
$ws.Range("G9").Value = 1
$ws.Range("E10").Value = 4
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 1
$ws.Range("E12").Value = 4
$ws.Range("G12").Value = 0.33
$ws.Range("E13").Value = 5
$ws.Range("G13").Value = 0.25

# Update the view: scroll position and selection
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 2
$ws.Range("J18").Select()

# Update page setup scale (keep "fit to height = 0 pages" in effect)
$ws.PageSetup.Zoom = 67
$ws.PageSetup.FitToPagesTall = 0
